# Update the "provenance" sheet: the python source file reference changed
# from tools_xl.py to its compiled counterpart tools_xl.pyc, and the
# generation timestamp moved forward a few minutes.
$wb = $excel.ActiveWorkbook

$wsProv = $wb.Worksheets.Item("provenance")
$wsProv.Range("B4").Value = "tools_xl.pyc"
$wsProv.Range("B12").Value = 43434.86570732237

# "08-BC" chapter sheet: start looping over chapter sections by adding the
# short section label "08-BC" a couple of rows below the existing title.
$wsBC = $wb.Worksheets.Item("08-BC")
$wsBC.Range("A3").Value = "08-BC"
